{"js": "// Apply the abstract.docx revision:\n//  1. Trim the author's first name, leaving \"  Winderbaum\".\n//  2. Add a soft hyphen between \"increasing\" and \"over\" so the line can break there.\n//  3. Expand the literature-review sentence with the list of challenges it uncovers.\n//  4. Insert a new sentence about aligning MathsTrack content with the senior\n//     high-school curricula, ahead of the existing alignment/recommendation sentence.\n//  5. Rewrite the closing sentence of the abstract.\n\nasync function replaceOnce(context, oldText, newText) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. \"Lyron Winderbaum\" -> \"  Winderbaum\" (two runs of \" \" + \" Winderbaum\" in the\n//    canonical XML; net visible text keeps a leading space and drops \"Lyron\").\nawait replaceOnce(context, \"Lyron Winderbaum\", \"  Winderbaum\");\n\n// 2. Insert a soft hyphen (U+00AD) between \"increasing \" and \"over\".\nawait replaceOnce(\n  context,\n  \"steadily increasing over the past several decades\",\n  \"steadily increasing \\u00ADover the past several decades\"\n);\n\n// 3. Literature review now uncovers specific challenges and promising approaches.\nawait replaceOnce(\n  context,\n  \"First, a literature review focused on challenges faced by mathematics bridging students and approaches to overcoming these challenges.\",\n  \"First, a literature review uncovers some of the key challenges faced by mathematics bridging students: maths anxiety, the secondary-tertiary education transition, negative affect towards maths, and self-efficacy. The literature also provides insight into approaches that show promise in overcoming these challenges.\"\n);\n\n// 4. Insert the new \"Alignment of the content of MathsTrack...\" sentence before the\n//    existing \"Alignment of these curricula is discussed...\" sentence.\nawait replaceOnce(\n  context,\n  \"(in particular MathsTrack). Alignment of these curricula is discussed, and recommendations made that could bring MathsTrack into closer alignment with the current high school curricula, within the context of the existing\",\n  \"(in particular MathsTrack). Alignment of the content of MathsTrack to the Australian senior high school curricula is desirable in order to support students in their future success as many of the common pathways they will pursue (entry level university maths courses, for example), are designed based on the assumption students entering them are coming from having completed high school in Australia recently.  Alignment of these curricula is discussed, and recommendations made that could bring MathsTrack into closer alignment with the current high school curricula, within the context of the existing\"\n);\n\n// 5. Rewrite the closing sentence of the abstract.\nawait replaceOnce(\n  context,\n  \"These recommendations will hopefully offer avenues for continued improvement for the University of Adelaide\\u2019s bridging courses, and beyond that provide some broader context for mathematics bridging in general.\",\n  \"The recommendations made from the synthesis of these two avenues of research will hopefully offer actionable pathways for the continued improvement for the University of Adelaide\\u2019s bridging courses, and beyond that provide some broader context of the concerns that need to be addressed in mathematics bridging more broadly.\"\n);\n", "ps1": "# Apply the abstract.docx revision:\n#  1. Trim the author's first name, leaving \"  Winderbaum\".\n#  2. Add a soft hyphen between \"increasing\" and \"over\" so the line can break there.\n#  3. Expand the literature-review sentence with the list of challenges it uncovers.\n#  4. Insert a new sentence about aligning MathsTrack content with the senior\n#     high-school curricula, ahead of the existing alignment/recommendation sentence.\n#  5. Rewrite the closing sentence of the abstract.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once($doc, $findText, $replaceText) {\n    $range = $doc.Content\n    $found = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n# 1. \"Lyron Winderbaum\" -> \"  Winderbaum\" (two runs of \" \" + \" Winderbaum\" in the\n#    canonical XML; net visible text keeps a leading space and drops \"Lyron\").\nReplace-Once $d \"Lyron Winderbaum\" \"  Winderbaum\"\n\n# 2. Insert a soft hyphen (U+00AD) between \"increasing \" and \"over\".\n$softHyphen = [char]0x00AD\n$find2 = \"steadily increasing over the past several decades\"\n$replace2 = \"steadily increasing \" + $softHyphen + \"over the past several decades\"\nReplace-Once $d $find2 $replace2\n\n# 3. Literature review now uncovers specific challenges and promising approaches.\nReplace-Once $d \"First, a literature review focused on challenges faced by mathematics bridging students and approaches to overcoming these challenges.\" \"First, a literature review uncovers some of the key challenges faced by mathematics bridging students: maths anxiety, the secondary-tertiary education transition, negative affect towards maths, and self-efficacy. The literature also provides insight into approaches that show promise in overcoming these challenges.\"\n\n# 4. Insert the new \"Alignment of the content of MathsTrack...\" sentence before the\n#    existing \"Alignment of these curricula is discussed...\" sentence.\nReplace-Once $d \"(in particular MathsTrack). Alignment of these curricula is discussed, and recommendations made that could bring MathsTrack into closer alignment with the current high school curricula, within the context of the existing\" \"(in particular MathsTrack). Alignment of the content of MathsTrack to the Australian senior high school curricula is desirable in order to support students in their future success as many of the common pathways they will pursue (entry level university maths courses, for example), are designed based on the assumption students entering them are coming from having completed high school in Australia recently.  Alignment of these curricula is discussed, and recommendations made that could bring MathsTrack into closer alignment with the current high school curricula, within the context of the existing\"\n\n# 5. Rewrite the closing sentence of the abstract.\n$rightQuote = [char]0x2019\n$find5 = \"These recommendations will hopefully offer avenues for continued improvement for the University of Adelaide\" + $rightQuote + \"s bridging courses, and beyond that provide some broader context for mathematics bridging in general.\"\n$replace5 = \"The recommendations made from the synthesis of these two avenues of research will hopefully offer actionable pathways for the continued improvement for the University of Adelaide\" + $rightQuote + \"s bridging courses, and beyond that provide some broader context of the concerns that need to be addressed in mathematics bridging more broadly.\"\nReplace-Once $d $find5 $replace5\n"}
